$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Replace("100 Australian speces", "100 Australian species", 1)
